# The edit removes the 5 "per-Million" / "per-100rb" metric columns
# (New_Cases_per_Million, Total_Cases_per_Million, New_Death_per_Million,
# Total_Deaths_per_Million, Total_Deaths_per_100rb) that used to live in
# columns L:P. Deleting those columns shifts the remaining
# Case_Fatality_Rate / Case_Recovered_Rate / Growth_Factor_of_New_Cases /
# Growth_Factor_of_New_Deaths columns (previously Q:T) left into L:O,
# which is exactly what the commit ("linear regression & svm" - dropping
# unused per-capita columns) shows in the target workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Columns("L:P").Delete()
